# Update odds values in Jogos_do_Dia_Betfair_Back_Lay_2026-01-06.xlsx (Sheet1)
# Columns use 1-based indices: F=6, G=7, H=8, I=9, J=10, K=11, L=12, M=13, N=14,
# O=15, P=16, Q=17, R=18, S=19, T=20, U=21, V=22, W=23, X=24, Y=25, Z=26,
# AA=27, AB=28, AC=29, AD=30, AE=31, AF=32, AG=33, AH=34, AI=35, AJ=36,
# AK=37, AL=38, AM=39, AN=40, AO=41

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Melbourne City vs Brisbane Roar)
$ws.Cells.Item(2, 16).Value = 1.81
$ws.Cells.Item(2, 18).Value = 1.31
$ws.Cells.Item(2, 25).Value = 15
$ws.Cells.Item(2, 28).Value = 7.8
$ws.Cells.Item(2, 35).Value = 80

# Row 3 (Pisa vs Como)
$ws.Cells.Item(3, 9).Value = 1.81
$ws.Cells.Item(3, 16).Value = 1.79
$ws.Cells.Item(3, 21).Value = 1.89
$ws.Cells.Item(3, 41).Value = 14

# Row 4 (Lecce vs Roma)
$ws.Cells.Item(4, 8).Value = 1.81
$ws.Cells.Item(4, 9).Value = 1.82
$ws.Cells.Item(4, 14).Value = 2.98
$ws.Cells.Item(4, 39).Value = 200

# Row 5 (Sassuolo vs Juventus)
$ws.Cells.Item(5, 10).Value = 3.95
$ws.Cells.Item(5, 11).Value = 4
$ws.Cells.Item(5, 13).Value = 1.07
$ws.Cells.Item(5, 16).Value = 1.91
$ws.Cells.Item(5, 18).Value = 1.34
$ws.Cells.Item(5, 20).Value = 1.97

# Row 6 (Livingston vs St Mirren)
$ws.Cells.Item(6, 6).Value = 2.78
$ws.Cells.Item(6, 7).Value = 2.88
$ws.Cells.Item(6, 8).Value = 2.66
$ws.Cells.Item(6, 9).Value = 2.74
$ws.Cells.Item(6, 10).Value = 3.6
$ws.Cells.Item(6, 11).Value = 3.65
$ws.Cells.Item(6, 12).Value = 1.39
$ws.Cells.Item(6, 14).Value = 3.5
$ws.Cells.Item(6, 15).Value = 1.35
$ws.Cells.Item(6, 16).Value = 1.84
$ws.Cells.Item(6, 17).Value = 2.06
$ws.Cells.Item(6, 18).Value = 1.32
$ws.Cells.Item(6, 19).Value = 3.75
$ws.Cells.Item(6, 20).Value = 1.83
$ws.Cells.Item(6, 21).Value = 2.12
$ws.Cells.Item(6, 22).Value = 1.58
$ws.Cells.Item(6, 23).Value = 1.53
$ws.Cells.Item(6, 24).Value = 14
$ws.Cells.Item(6, 25).Value = 11
$ws.Cells.Item(6, 26).Value = 16.5
$ws.Cells.Item(6, 27).Value = 40
$ws.Cells.Item(6, 28).Value = 11.5
$ws.Cells.Item(6, 29).Value = 7.8
$ws.Cells.Item(6, 30).Value = 12.5
$ws.Cells.Item(6, 31).Value = 30
$ws.Cells.Item(6, 32).Value = 18.5
$ws.Cells.Item(6, 33).Value = 12.5
$ws.Cells.Item(6, 34).Value = 19.5
$ws.Cells.Item(6, 36).Value = 46
$ws.Cells.Item(6, 37).Value = 32
$ws.Cells.Item(6, 40).Value = 32
$ws.Cells.Item(6, 41).Value = 26

# Row 7 (West Ham vs Nottm Forest)
$ws.Cells.Item(7, 12).Value = 1.41
$ws.Cells.Item(7, 16).Value = 1.96
$ws.Cells.Item(7, 23).Value = 1.42
$ws.Cells.Item(7, 28).Value = 13
$ws.Cells.Item(7, 32).Value = 22
$ws.Cells.Item(7, 41).Value = 19

# Row 8 (Rangers vs Aberdeen)
$ws.Cells.Item(8, 6).Value = 1.63
$ws.Cells.Item(8, 7).Value = 1.64
$ws.Cells.Item(8, 8).Value = 6.6
$ws.Cells.Item(8, 9).Value = 7.4
$ws.Cells.Item(8, 10).Value = 4.1
$ws.Cells.Item(8, 11).Value = 4.2
$ws.Cells.Item(8, 14).Value = 4.8
$ws.Cells.Item(8, 15).Value = 1.25
$ws.Cells.Item(8, 16).Value = 2.26
$ws.Cells.Item(8, 17).Value = 1.74
$ws.Cells.Item(8, 18).Value = 1.53
$ws.Cells.Item(8, 19).Value = 2.72
$ws.Cells.Item(8, 20).Value = 1.79
$ws.Cells.Item(8, 21).Value = 2.1
$ws.Cells.Item(8, 22).Value = 1.16
$ws.Cells.Item(8, 23).Value = 2.54
$ws.Cells.Item(8, 24).Value = 21
$ws.Cells.Item(8, 25).Value = 27
$ws.Cells.Item(8, 26).Value = 60
$ws.Cells.Item(8, 27).Value = 190
$ws.Cells.Item(8, 28).Value = 10.5
$ws.Cells.Item(8, 29).Value = 9.800000000000001
$ws.Cells.Item(8, 30).Value = 26
$ws.Cells.Item(8, 31).Value = 90
$ws.Cells.Item(8, 33).Value = 10
$ws.Cells.Item(8, 34).Value = 19
$ws.Cells.Item(8, 35).Value = 85
$ws.Cells.Item(8, 37).Value = 15.5
$ws.Cells.Item(8, 38).Value = 32
$ws.Cells.Item(8, 39).Value = 90
$ws.Cells.Item(8, 40).Value = 7.6
$ws.Cells.Item(8, 41).Value = 100
